$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data for rows 2-11 (columns A-G):
# employee_id, employee_name, department, absence_reason, absence_duration, absence_date, salary
$data = @(
    @(60958, "Dra. Maria Luísa Correia", "P&D", "Consulta medica", 2, 45100, 4705.31),
    @(48719, "Caroline Costela", "TI", "Problemas pessoais", 5, 45086, 6666.29),
    @(69099, "Yasmin Fonseca", "P&D", "Consulta medica", 3, 45083, 6595.34),
    @(70228, "Luiz Gustavo Novaes", "Financeiro", "Outros", 1, 45093, 4163.93),
    @(377, "Mirella da Rosa", "Operacoes", "Problemas pessoais", 2, 45084, 2005.19),
    @(51355, "Valentim Machado", "P&D", "Consulta medica", 7, 45106, 3735.47),
    @(28171, "Eduardo Nascimento", "Vendas", "Doenca", 7, 45083, 5365.11),
    @(11456, "Sr. Diego Aragão", "Financeiro", "Outros", 2, 45080, 3429.58),
    @(33119, "Maria Vitória da Cunha", "Financeiro", "Outros", 5, 45105, 6550.48),
    @(1315, "João Lucas Cirino", "Juridico", "Consulta medica", 8, 45083, 2800.15)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $r++
}
